# RPA datasets push 2024-08-01
# Insert a new IPO record ("피앤에스미캐닉스") as the first data row of the
# table, shifting every existing row down by one. All other rows/data are
# otherwise unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push all existing data rows down by one, landing the new record in row 2
# (right after the header row).
$ws.Rows(2).Insert()
$ws.Rows(2).ClearFormats()

# A2 / O2 / P2 hold date-like text ("2024-07-31", "2024-07-22", "2024-07-25").
# Force a text number format first so Excel stores them as plain strings
# instead of silently converting them to date serial numbers.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("O2").NumberFormat = "@"
$ws.Range("P2").NumberFormat = "@"

$ws.Range("A2").Value = "2024-07-31"
$ws.Range("B2").Value = "피앤에스미캐닉스"
$ws.Range("C2").Value = "코스닥"
$ws.Range("D2").Value = 297
$ws.Range("E2").Value = "키움"
$ws.Range("F2").Value = 297
$ws.Range("G2").Value = "-"
$ws.Range("H2").Value = "-"
$ws.Range("I2").Value = "-"
$ws.Range("J2").Value = "-"
$ws.Range("K2").Value = "대표"
$ws.Range("L2").Value = "-"
$ws.Range("M2").Value = 22000
$ws.Range("N2").Value = 100
$ws.Range("O2").Value = "2024-07-22"
$ws.Range("P2").Value = "2024-07-25"
$ws.Range("Q2").Value = 1012500

# Drop the temporary text format so the new row matches the plain
# (unstyled) formatting used by the rest of the data rows.
$ws.Rows(2).ClearFormats()
